# Update the cryptocurrency table with refreshed prices / volumes.
# Rows for several coins were also reordered / replaced (B37-B40, B44-B45, B51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force Text number format so numeric-looking strings (e.g. "46.41")
    # are not silently coerced into Excel numbers, matching the source data
    # which stores these as inline/shared strings.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell "D2" "30.851.23"
Set-TextCell "E2" "  +2.28%  "
Set-TextCell "D3" "2.121.62"
Set-TextCell "E3" "  +10.35%  "
Set-TextCell "D5" "335.46"
Set-TextCell "E5" "  +4.89%  "
Set-TextCell "D6" "0.9998"
Set-TextCell "E6" "  -0.17%  "
Set-TextCell "D7" "0.5382"
Set-TextCell "E7" "  +6.04%  "
Set-TextCell "D8" "0.4432"
Set-TextCell "E8" "  +8.70%  "
Set-TextCell "D9" "0.09098"
Set-TextCell "E9" "  +9.04%  "
Set-TextCell "D10" "46.41"
Set-TextCell "E10" "  +10.39%  "
Set-TextCell "D11" "1.186"
Set-TextCell "E11" "  +5.94%  "
Set-TextCell "D12" "25.44"
Set-TextCell "E12" "  +5.78%  "
Set-TextCell "D13" "2.118.72"
Set-TextCell "E13" "  +10.19%  "
Set-TextCell "D14" "6.799"
Set-TextCell "E14" "  +5.71%  "
Set-TextCell "D15" "7.868"
Set-TextCell "E15" "  +8.42%  "
Set-TextCell "D16" "98.44"
Set-TextCell "E16" "  +6.14%  "
Set-TextCell "E17" "  +4.36%  "
Set-TextCell "D19" "0.06657"
Set-TextCell "E19" "  +2.24%  "
Set-TextCell "E20" "  +4.16%  "
Set-TextCell "D21" "6.413"
Set-TextCell "E21" "  +7.59%  "
Set-TextCell "E22" "  -0.13%  "
Set-TextCell "D23" "30.963.52"
Set-TextCell "E23" "  +2.63%  "
Set-TextCell "D24" "12.24"
Set-TextCell "E24" "  +7.62%  "
Set-TextCell "D25" "2.373.21"
Set-TextCell "E25" "  +10.75%  "
Set-TextCell "E26" "  +3.61%  "
Set-TextCell "D27" "22.97"
Set-TextCell "E27" "  +4.48%  "
Set-TextCell "D28" "2.615"
Set-TextCell "E28" "  +15.34%  "
Set-TextCell "D29" "163.95"
Set-TextCell "D30" "134.19"
Set-TextCell "E30" "  +4.02%  "
Set-TextCell "D31" "1.177"
Set-TextCell "E31" "  +3.49%  "
Set-TextCell "D32" "0.1085"
Set-TextCell "E32" "  +3.67%  "
Set-TextCell "D33" "6.305"
Set-TextCell "E33" "  +5.75%  "
Set-TextCell "D34" "4.015"
Set-TextCell "E34" "  +5.64%  "
Set-TextCell "D35" "1.553"
Set-TextCell "E35" "  +27.56%  "
Set-TextCell "D36" "0.02621"
Set-TextCell "E36" "  +6.87%  "
Set-TextCell "B37" "Aptos"
Set-TextCell "C37" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D37" "13.38"
Set-TextCell "E37" "  +16.80%  "
Set-TextCell "B38" "InternetComputer(DFINITY)"
Set-TextCell "C38" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D38" "5.618"
Set-TextCell "E38" "  +5.56%  "
Set-TextCell "B39" "FraxShare"
Set-TextCell "C39" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D39" "9.656"
Set-TextCell "E39" "  +12.33%  "
Set-TextCell "B40" "Hedera"
Set-TextCell "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D40" "0.06776"
Set-TextCell "E40" "  +5.12%  "
Set-TextCell "D41" "0.2284"
Set-TextCell "E41" "  +6.42%  "
Set-TextCell "D42" "0.6886"
Set-TextCell "E42" "  +5.82%  "
Set-TextCell "D43" "1.262"
Set-TextCell "E43" "  +4.19%  "
Set-TextCell "B44" "EnergySwap"
Set-TextCell "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D44" "14.20"
Set-TextCell "E44" "  +5.78%  "
Set-TextCell "B45" "Decentraland"
Set-TextCell "C45" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell "D45" "0.6470"
Set-TextCell "E45" "  +6.79%  "
Set-TextCell "D46" "0.9994"
Set-TextCell "E46" "  -0.14%  "
Set-TextCell "D47" "2.262"
Set-TextCell "E47" "  +3.56%  "
Set-TextCell "D48" "3.690"
Set-TextCell "E48" "  +1.76%  "
Set-TextCell "D49" "1.294"
Set-TextCell "D50" "83.65"
Set-TextCell "E50" "  +7.07%  "
Set-TextCell "B51" "WEMIXTOKEN"
Set-TextCell "C51" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D51" "1.167"
Set-TextCell "E51" "  +3.02%  "
